# Commit: "added hybrid changes to sitePreparationCost"
#
# Observed edits in the workbook:
#   1. The "Collection mode" input (column AT, rows 2 & 3 on the
#      "Project list" sheet) changes from "manual" to "auto".
#   2. The selected cell on the "Project list" sheet moves to AU40.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project list")

# 1) Collection mode: "manual" -> "auto" for both rows that use it.
$ws.Range("AT2").Value = "auto"
$ws.Range("AT3").Value = "auto"

# 2) Update the active selection on the sheet.
$ws.Range("AU40").Select()
